# Fill in previously-missing daily spot price columns (DF:DK, i.e.
# 01-oct. through 06-oct.) for every hour row (2-25) on the "Prix Spot"
# sheet. These cells currently hold an inline "-" placeholder string;
# replace them with the actual numeric price values now available.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

$data = @{
  2  = @(79.53, 69.26000000000001, 32.65, 9.85, 4.31, 26.96)
  3  = @(79.98, 66.84999999999999, 36.09, 7.23, 1.63, 21.61)
  4  = @(66.95, 68.91, 34.49, 2.9, 0.58, 16.38)
  5  = @(37.05, 42.07, 25.68, 0.79, 0.11, 15.9)
  6  = @(28.49, 25.05, 21.94, -0, -0.01, 15.92)
  7  = @(44.03, 39.44, 19.01, 0, -0.01, 20.9)
  8  = @(31.72, 55.25, 28.41, 0, -0.01, 18.28)
  9  = @(40.07, 70.3, 28.13, -0.01, -0.03, 46.92)
  10 = @(73.56, 75.23, 65.41, 0.14, 1.36, 91.58)
  11 = @(106.73, 110.24, 90.14, 1.85, 3.78, 115.81)
  12 = @(120.44, 122.97, 103.28, 6.6, 2.22, 124.76)
  13 = @(110.8, 109.62, 99.3, 11.46, 3.14, 118.02)
  14 = @(80.23, 73.89, 76.89, 8.710000000000001, 0.67, 97.73)
  15 = @(57.53, 48.64, 58.94, 2.21, -0.03, 62.72)
  16 = @(36.52, 31.83, 47.38, 0.04, -0.49, 59.93)
  17 = @(23.92, 22.1, 41.79, -0.01, -0.99, 58.24)
  18 = @(18.39, 19.43, 41.44, -0.01, -1.84, 46.95)
  19 = @(13.43, 16.2, 36.8, -0.03, -2.48, 37.62)
  20 = @(16.2, 16.25, 38.24, -0.01, -2.29, 36.23)
  21 = @(19.89, 16.99, 39.56, -0, -0.78, 41.36)
  22 = @(36.97, 31.89, 44.1, -0.01, -0, 49.19)
  23 = @(43.02, 38.28, 49.54, -0.02, 3.03, 63.01)
  24 = @(66.23, 55.38, 52.51, 1.64, 16.58, 82.58)
  25 = @(101.11, 94.09, 57.36, 16.57, 34.34, 117.46)
}

$cols = @("DF", "DG", "DH", "DI", "DJ", "DK")

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])$row").Value = $vals[$i]
  }
}
